$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.581.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.603.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.129"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.054.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60.620.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000140"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.611.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "354.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("E20").Value = "  +4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.97%  "
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.431"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.719.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.948"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.35%  "
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.844"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.629"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.982.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
